$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Tv3 C4"
$ws.Range("B2").Value = 2300
$ws.Range("C2").Value = 280

$ws.Range("A3").Value = "Ck1 White Diamond"
$ws.Range("B3").Value = 2200
$ws.Range("C3").Value = 307

$ws.Range("B3").Select()
